$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2355
$ws.Range("B3").Value = 2041
$ws.Range("B9").Value = 745
$ws.Range("B16").Value = 708
